# Auto-generated Excel COM-interop script
# Updates cached market-price / profit values in the Kujata_Profits workbook
# (each worksheet tab = a Disciple of the Hand/Land class, e.g. ALC, ARM, BSM, ...)
# following a scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")

$updates = [ordered]@{
    "H100" = 500
    "I100" = 500
    "K100" = 500
    "M100" = 41
    "H106" = 2575.875
    "I106" = 2600.6667
    "K106" = 2600.6667
    "M106" = -1969.6667
    "H13" = 1600
    "J13" = 1600
    "L13" = 1600
    "N13" = -1938
    "H137" = 1136.4762
    "I137" = 1125.7222
    "J137" = 1201
    "K137" = 3377.1666
    "L137" = 3603
    "M137" = -827.1665999999996
    "N137" = -8703
    "H33" = 406.77777
    "I33" = 425.85715
    "J33" = 340
    "K33" = 425.85715
    "L33" = 340
    "M33" = -196.85715
    "N33" = -798
    "H64" = 3877.4285
    "I64" = 4660.6665
    "J64" = 3290
    "K64" = 4660.6665
    "L64" = 3290
    "M64" = -4412.6665
    "N64" = -3786
    "H67" = 3877.4285
    "I67" = 4660.6665
    "J67" = 3290
    "K67" = 4660.6665
    "L67" = 3290
    "M67" = -3802.6665
    "N67" = -5006
    "H74" = 8366.944
    "I74" = 10226.333
    "J74" = 4648.1665
    "K74" = 10226.333
    "L74" = 4648.1665
    "M74" = -9290.333000000001
    "N74" = -6520.1665
    "H77" = 8366.944
    "I77" = 10226.333
    "J77" = 4648.1665
    "K77" = 51131.665
    "L77" = 23240.8325
    "M77" = -46451.665
    "N77" = -32600.8325
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")

$updates = [ordered]@{
    "H132" = 2723.8845
    "I132" = 2351.9412
    "J132" = 3426.4443
    "K132" = 7055.823600000001
    "L132" = 10279.3329
    "M132" = -4525.823600000001
    "N132" = -15339.3329
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")

$updates = [ordered]@{
    "H105" = 66668640
    "I105" = 83335170
    "J105" = 2540
    "K105" = 83335170
    "L105" = 2540
    "M105" = -83333423
    "N105" = -6034
    "H86" = 5389.1113
    "I86" = 5467.2666
    "J86" = 4998.3335
    "K86" = 5467.2666
    "L86" = 4998.3335
    "M86" = -4344.2666
    "N86" = -7244.3335
    "H89" = 5389.1113
    "I89" = 5467.2666
    "J89" = 4998.3335
    "K89" = 27336.333
    "L89" = 24991.6675
    "M89" = -21720.333
    "N89" = -36223.6675
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")

$updates = [ordered]@{
    "H105" = 479.14285
    "I105" = 425.66666
    "J105" = 800
    "K105" = 425.66666
    "L105" = 800
    "M105" = 1321.33334
    "N105" = -4294
    "H132" = 2550.5386
    "I132" = 1351.2858
    "K132" = 4053.8574
    "M132" = -1523.8574
    "H134" = 55557256
    "I134" = 83334630
    "J134" = 2500
    "K134" = 250003890
    "L134" = 7500
    "M134" = -250001355
    "N134" = -12570
    "H136" = 1314.1818
    "I136" = 1120
    "J136" = 1476
    "K136" = 3360
    "L136" = 4428
    "M136" = -810
    "N136" = -9528
    "H22" = 474.5
    "I22" = 466.33334
    "J22" = 499
    "K22" = 466.33334
    "L22" = 499
    "M22" = -116.33334
    "N22" = -1199
    "H5" = 1009.3333
    "I5" = 1399
    "J5" = 230
    "K5" = 1399
    "L5" = 230
    "M5" = -1287
    "N5" = -454
    "H58" = 1314.1818
    "I58" = 1120
    "J58" = 1476
    "K58" = 1120
    "L58" = 1476
    "M58" = -917
    "N58" = -1882
    "H7" = 375.16666
    "I7" = 375.16666
    "J7" = 0
    "K7" = 375.16666
    "L7" = 0
    "M7" = -262.16666
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

$clears = @("N7")
foreach ($ref in $clears) {
    $ws.Range($ref).ClearContents()
}

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")

$updates = [ordered]@{
    "H11" = 198
    "I11" = 198
    "J11" = 0
    "K11" = 594
    "L11" = 0
    "M11" = -454
    "H122" = 838.3
    "I122" = 772.1667
    "J122" = 937.5
    "K122" = 6949.5003
    "L122" = 8437.5
    "M122" = -4499.5003
    "N122" = -13337.5
    "H131" = 14707139
    "J131" = 1345.1451
    "L131" = 4035.4353
    "N131" = -14115.4353
    "H132" = 1222.5
    "I132" = 950
    "J132" = 1495
    "K132" = 8550
    "L132" = 13455
    "M132" = -6020
    "N132" = -18515
    "H14" = 586
    "I14" = 586
    "K14" = 1758
    "M14" = -1585
    "H68" = 1183.3334
    "I68" = 1100
    "K68" = 3300
    "M68" = -2489
    "H69" = 2112.611
    "I69" = 1149.75
    "J69" = 2387.7144
    "K69" = 3449.25
    "L69" = 7163.1432
    "M69" = -2638.25
    "N69" = -8785.143199999999
    "H71" = 1183.3334
    "I71" = 1100
    "K71" = 9900
    "M71" = -5844
    "H72" = 2112.611
    "I72" = 1149.75
    "J72" = 2387.7144
    "K72" = 10347.75
    "L72" = 21489.4296
    "M72" = -6291.75
    "N72" = -29601.4296
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

$clears = @("N11")
foreach ($ref in $clears) {
    $ws.Range($ref).ClearContents()
}

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")

$updates = [ordered]@{
    "H109" = 29247.75
    "J109" = 30000
    "L109" = 30000
    "N109" = -32080
    "H126" = 2306.625
    "I126" = 1790.8572
    "J126" = 2707.7778
    "K126" = 5372.571599999999
    "L126" = 8123.3334
    "M126" = -2902.571599999999
    "N126" = -13063.3334
    "H21" = 5000000
    "I21" = 5000000
    "K21" = 5000000
    "M21" = -4999827
    "H30" = 5000000
    "I30" = 5000000
    "K30" = 5000000
    "M30" = -4999895
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")

$updates = [ordered]@{
    "H104" = 0
    "J104" = 0
    "L104" = 0
    "H113" = 1259.4445
    "I113" = 1266.875
    "J113" = 1200
    "K113" = 1266.875
    "L113" = 1200
    "M113" = 903.125
    "N113" = -5540
    "H132" = 93163.91
    "I132" = 1333.3334
    "J132" = 203360.6
    "K132" = 4000.0002
    "L132" = 610081.8
    "M132" = -1470.0002
    "N132" = -615141.8
    "H136" = 1568.7273
    "I136" = 1282.0625
    "J136" = 2333.1667
    "K136" = 3846.1875
    "L136" = 6999.500100000001
    "M136" = -1296.1875
    "N136" = -12099.5001
    "H22" = 1573.25
    "I22" = 1433.3334
    "J22" = 1993
    "K22" = 1433.3334
    "L22" = 1993
    "M22" = -1138.3334
    "N22" = -2583
    "H27" = 1573.25
    "I27" = 1433.3334
    "J27" = 1993
    "K27" = 1433.3334
    "L27" = 1993
    "M27" = -1326.3334
    "N27" = -2207
    "H55" = 221.86667
    "J55" = 425.33334
    "L55" = 425.33334
    "N55" = -771.33334
    "H61" = 1259.4445
    "I61" = 1266.875
    "J61" = 1200
    "K61" = 1266.875
    "L61" = 1200
    "M61" = -1064.875
    "N61" = -1604
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

$clears = @("N104")
foreach ($ref in $clears) {
    $ws.Range($ref).ClearContents()
}

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")

$updates = [ordered]@{
    "H12" = 46838
    "J12" = 46838
    "L12" = 46838
    "N12" = -47122
    "H122" = 12383341
    "I122" = 12383341
    "J122" = 0
    "K122" = 37150023
    "L122" = 0
    "M122" = -37147573
    "H132" = 4353.2144
    "I132" = 4261.4
    "J132" = 4582.75
    "K132" = 12784.2
    "L132" = 13748.25
    "M132" = -10254.2
    "N132" = -18808.25
    "H136" = 1012.1429
    "I136" = 521.25
    "K136" = 1563.75
    "M136" = 986.25
    "H96" = 1822.8
    "I96" = 2962.5
    "J96" = 1063
    "K96" = 2962.5
    "L96" = 1063
    "M96" = -1589.5
    "N96" = -3809
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

$clears = @("N122")
foreach ($ref in $clears) {
    $ws.Range($ref).ClearContents()
}

